$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 0.65
$ws.Range("B16").Value = 1657238.437593
$ws.Range("C16").Value = 20357.699105
$ws.Range("D16").Value = 1636880.738488
$ws.Range("E16").Value = 12974.03040366667
$ws.Range("F16").Value = 748182.757149
$ws.Range("G16").Value = 14232.5759
$ws.Range("H16").Value = 733950.18125
$ws.Range("I16").Value = 14840.43002133333
